$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("weeknr 48")

# Row 17: fill in the StartScene log entry (begin/end time, activity text)
$ws.Range("C17").Value = 0.57291666666666663
$ws.Range("D17").Value = 0.59375
$ws.Range("F17").Value = "Scene waardes ingevuld. Moeten nog wat fouten uitgehaald worden."

# Row grows to two lines once the activity text is filled in
$ws.Rows.Item(17).RowHeight = 30

# Move the active selection from F15 to F17 (where the new entry was typed)
$ws.Range("F17").Select()
